$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 4 (pushes the old "total" row 4 down to row 6)
$ws.Rows.Item(4).Insert()
$ws.Rows.Item(4).Insert()

# --- Row 2: update contract/id/name, and rework the rappel-related columns ---
$ws.Range("A2").Value = "001/RRR/AV1"
$ws.Range("C2").Value = "B219321"
$ws.Range("D2").Value = "JEMAA HORMI"
$ws.Range("H2").Value = "--"
$ws.Range("J2").Value = "--"
$ws.Range("L2").Value = 8000
$ws.Range("N2").Value = "--"
$ws.Range("O2").Value = 8000

# --- Row 3: update contract number, and rework the rappel-related columns ---
$ws.Range("A3").Value = "001/RRR/AV1"
$ws.Range("H3").Value = "--"
$ws.Range("J3").Value = "--"
$ws.Range("L3").Value = 8000
$ws.Range("N3").Value = "--"
$ws.Range("O3").Value = 8000

# --- New row 4: JEMAA HORMI monthly entry (non-rappel) ---
$ws.Range("A4").Value = "001/RRR/AV1"
$ws.Range("B4").Value = "Direction régionale"
$ws.Range("C4").Value = "B219321"
$ws.Range("D4").Value = "JEMAA HORMI"
$ws.Range("E4").Value = "non"
$ws.Range("F4").Value = "mensuelle"
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 1000
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = "--"
$ws.Range("O4").Value = 1000

# --- New row 5: ZERNAKH ABDELLAH monthly entry (non-rappel) ---
$ws.Range("A5").Value = "001/RRR/AV1"
$ws.Range("B5").Value = "Direction régionale"
$ws.Range("C5").Value = "IB19558"
$ws.Range("D5").Value = "ZERNAKH ABDELLAH"
$ws.Range("E5").Value = "non"
$ws.Range("F5").Value = "mensuelle"
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 1000
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 0
$ws.Range("N5").Value = "--"
$ws.Range("O5").Value = 1000

# --- Row 6 (former totals row, shifted down by the inserts): update totals ---
$ws.Range("L6").Value = 16000
$ws.Range("O6").Value = 18000
